$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.156.88"
$ws.Range("E2").Value = "  -2.72%  "
$ws.Range("D3").Value = "'1.871.10"
$ws.Range("E3").Value = "  -2.05%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'307.13"
$ws.Range("E5").Value = "  -1.79%  "
$ws.Range("D6").Value = "'0.9997"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "'0.5129"
$ws.Range("E7").Value = "  +2.51%  "
$ws.Range("D8").Value = "'0.3753"
$ws.Range("E8").Value = "  -1.34%  "
$ws.Range("D9").Value = "'0.07139"
$ws.Range("E9").Value = "  -2.00%  "
$ws.Range("D10").Value = "'0.8876"
$ws.Range("E10").Value = "  -2.39%  "
$ws.Range("D11").Value = "'20.63"
$ws.Range("E11").Value = "  -3.24%  "
$ws.Range("D12").Value = "'1.884.73"
$ws.Range("E12").Value = "  -1.40%  "
$ws.Range("D13").Value = "'0.07544"
$ws.Range("E13").Value = "  -1.29%  "
$ws.Range("D14").Value = "'5.324"
$ws.Range("E14").Value = "  -2.74%  "
$ws.Range("D15").Value = "'89.17"
$ws.Range("E15").Value = "  -3.75%  "
$ws.Range("D16").Value = "'0.9997"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("E17").Value = "  -2.93%  "
$ws.Range("E18").Value = "  -3.77%  "
$ws.Range("D19").Value = "'0.9992"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").Value = "'27.227.28"
$ws.Range("E20").Value = "  -2.62%  "
$ws.Range("D21").Value = "'5.048"
$ws.Range("E21").Value = "  -2.40%  "
$ws.Range("D22").Value = "'2.113.78"
$ws.Range("E22").Value = "  -2.47%  "
$ws.Range("E23").Value = "  -2.99%  "
$ws.Range("D24").Value = "'6.479"
$ws.Range("E24").Value = "  -1.94%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "'1.854"
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'149.96"
$ws.Range("E26").Value = "  -1.88%  "
$ws.Range("D27").Value = "'17.95"
$ws.Range("E27").Value = "  -2.42%  "
$ws.Range("D28").Value = "'2.096"
$ws.Range("E28").Value = "  -6.03%  "
$ws.Range("D29").Value = "'112.85"
$ws.Range("E29").Value = "  -1.95%  "
$ws.Range("D30").Value = "'4.714"
$ws.Range("E30").Value = "  -3.83%  "
$ws.Range("D31").Value = "'4.677"
$ws.Range("E31").Value = "  -2.76%  "
$ws.Range("D32").Value = "'0.09025"
$ws.Range("E32").Value = "  +0.47%  "
$ws.Range("E33").Value = "  -2.89%  "
$ws.Range("D34").Value = "'3.085"
$ws.Range("E34").Value = "  -3.57%  "
$ws.Range("D35").Value = "'1.158"
$ws.Range("E35").Value = "  -5.99%  "
$ws.Range("D36").Value = "'0.7340"
$ws.Range("E36").Value = "  -6.87%  "
$ws.Range("D37").Value = "'0.02049"
$ws.Range("E37").Value = "  -1.84%  "
$ws.Range("D38").Value = "'2.509"
$ws.Range("E38").Value = "  -5.13%  "
$ws.Range("D39").Value = "'3.059"
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("D40").Value = "'1.080"
$ws.Range("E40").Value = "  -0.93%  "
$ws.Range("D41").Value = "'0.5351"
$ws.Range("E41").Value = "  -3.36%  "
$ws.Range("D42").Value = "'6.579"
$ws.Range("E42").Value = "  -3.46%  "
$ws.Range("D43").Value = "'117.39"
$ws.Range("E43").Value = "  +3.43%  "
$ws.Range("E44").Value = "  -2.29%  "
$ws.Range("D45").Value = "'0.1476"
$ws.Range("E45").Value = "  -2.80%  "
$ws.Range("D46").Value = "'0.4635"
$ws.Range("E46").Value = "  -4.02%  "
$ws.Range("D47").Value = "'0.9993"
$ws.Range("D48").Value = "'10.01"
$ws.Range("D49").Value = "'1.570"
$ws.Range("E49").Value = "  -4.16%  "
$ws.Range("D50").Value = "'64.40"
$ws.Range("E50").Value = "  -4.35%  "
$ws.Range("D51").Value = "'36.50"
$ws.Range("E51").Value = "  -1.44%  "
